$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 560, shifting the existing rows 560-604 down to 561-605.
$ws.Rows(560).Insert()

# Populate the newly inserted row 560 with the new weekly data point.
$ws.Range("A560").Value = 10
$ws.Range("B560").Value = "Vega Modelo de Temuco"
$ws.Range("C560").Value = "La Araucanía"
$ws.Range("D560").Value = 45013
$ws.Range("E560").Value = 9
$ws.Range("F560").Value = 100112024
$ws.Range("G560").Value = "Choclo"
$ws.Range("H560").Value = "Dulce o Americano"
$ws.Range("I560").Value = "Primera"
$ws.Range("J560").Value = 5000
$ws.Range("K560").Value = 280
$ws.Range("L560").Value = 280
$ws.Range("M560").Value = 280
$ws.Range("N560").Value = '$/unidad'
$ws.Range("O560").Value = "Región de La Araucanía"
$ws.Range("P560").Value = 280
$ws.Range("Q560").Value = 1
$ws.Range("R560").Value = "Hortaliza"
